$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 83; this shifts the existing rows
# 83-151 down to 84-152 (carrying their data and formatting with them,
# which reproduces the "shift" seen throughout the diff).
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new data record.
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44904
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 100112021
$ws.Range("G83").Value = "Ají"
$ws.Range("H83").Value = "Inferno"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 50
$ws.Range("K83").Value = 14000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = 14400
$ws.Range("N83").Value = "$/caja 10 kilos"
$ws.Range("O83").Value = "Región de Arica y Parinacota"
$ws.Range("P83").Value = 1440
$ws.Range("Q83").Value = 10
$ws.Range("R83").Value = "Hortaliza"
